$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (new quarter columns) - copy formatting from Q1 (bold/centered/bordered header style)
$ws.Range("Q1").Copy()
$ws.Range("R1:T1").PasteSpecial(-4122)
$ws.Range("R1").Value = "31/12/2023"
$ws.Range("S1").Value = "31/03/2024"
$ws.Range("T1").Value = "30/06/2024"

# Data rows: R,S,T values for rows 2-80 (value rows get numbers, label-only rows get blanks)
$rowData = @{
    2 = @(2358531.072, 2633339.904, 3201582.08)
    3 = @(1549352.96, 1824141.952, 2310778.88)
    4 = @(465588.992, 493344.992, 404080)
    5 = @(264524.992, 212479.008, 420312.992)
    6 = @(483044, 437569.984, 56120)
    7 = @(138096, 416204.992, 798281.9840000001)
    8 = @(0, 0, 0)
    9 = @(56700, 54206, 56862)
    10 = @(0, 0, 0)
    11 = @(141399.008, 210336.992, 575121.9840000001)
    12 = @(141931.008, 143444, 145183.008)
    13 = @(0, 0, 0)
    14 = @(8372, 7199, 7887)
    15 = @(0, 0, 0)
    16 = @(0, 0, 0)
    17 = @(0, 0, 0)
    18 = @(0, 0, 0)
    19 = @(0, 0, 0)
    20 = @(0, 0, 0)
    21 = @(0, 0, 0)
    22 = @(1755, 1776, 1776)
    23 = @(663518.976, 661665.024, 741574.0159999999)
    24 = @(1973, 2313, 2270)
    25 = @(0, 0, 0)
    26 = @(2358531.072, 2633339.904, 3201582.08)
    27 = @(338088.992, 927038.0159999999, 1165126.016)
    28 = @(9899, 11498, 13291)
    29 = @(154324.992, 279784, 296876.992)
    30 = @(5787, 7751, 8172)
    31 = @(38533, 583078.0159999999, 712974.976)
    32 = @(0, 0, 0)
    33 = @(89386, 4584, 4778)
    34 = @(40159, 40343, 129033)
    35 = @(0, 0, 0)
    36 = @(0, 0, 0)
    37 = @(547934.976, 238788.992, 250758)
    38 = @(535056.992, 227904.992, 241684)
    39 = @(0, 0, 0)
    40 = @(12878, 10884, 9074)
    41 = @(0, 0, 0)
    42 = @(0, 0, 0)
    43 = @(0, 0, 0)
    44 = @(0, 0, 0)
    45 = @(0, 0, 0)
    46 = @(459878.016, 460204.992, 490726.016)
    47 = @(1012628.992, 1007307.968, 1294972.032)
    48 = @(429726.016, 429726.016, 719420.032)
    49 = @(1451, 1974, 2583)
    50 = @(0, 0, 0)
    51 = @(581452.032, 581452.032, 574577.9840000001)
    52 = @(0, -5844, -1609)
    53 = @(0, 0, 0)
    54 = @(0, 0, 0)
    55 = @(0, 0, 0)
    56 = @(0, 0, 0)
    59 = @(841084.096, 69102, 87572)
    60 = @(-695640.96, -76831, -56042)
    61 = @(145443.008, -7729, 31530)
    62 = @(-5598, -8335, -12291)
    63 = @(-10402, -8418, -4728)
    64 = @(-3586, -3820, -784)
    65 = @(-1139, 140, -673)
    66 = @(0, 0, 0)
    67 = @(79, 0, 0)
    68 = @(225, 33427, 3975)
    69 = @(10026, 76588, 35069)
    70 = @(-9801, -43161, -31094)
    74 = @(125022.008, 5265, 17029)
    75 = @(-721, 2826, -172)
    76 = @(91273, 0, 1257)
    79 = @(-23485, -13935, -13743)
    80 = @(191572.992, -5844, 4371)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 18).Value = $vals[0]
    $ws.Cells.Item($row, 19).Value = $vals[1]
    $ws.Cells.Item($row, 20).Value = $vals[2]
}

# Blank label-only rows: create empty (contentless) cells mirroring the existing blank
# cells already present in columns B-Q for these rows. A no-op format touch (LineStyle
# set to "no border", matching the existing formatting) makes Excel materialize the
# empty cell without giving it a value.
$blankRows = @(57, 58, 71, 72, 73, 77, 78)
foreach ($row in $blankRows) {
    $ws.Range("R" + $row + ":T" + $row).Borders.LineStyle = -4142
}

Write-Host "Done"